{"js": "// Insert three new bullet paragraphs right after the\n// \"GIS & Geospatial Analysis Consulting\" paragraph (under the Siege\n// Analytics / PARTNER heading), before the existing \"\u2022 Lead comprehensive\n// research initiatives...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"GIS & Geospatial Analysis Consulting\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === targetText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find paragraph: \" + targetText);\n}\n\nconst newLines = [\n  \"\\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n  \"\\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n  \"\\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n];\n\n// Insert the three paragraphs in order, each one right after the previous,\n// so they stay in the requested order immediately following the anchor.\nlet insertAfter = anchor;\nfor (const line of newLines) {\n  insertAfter = insertAfter.insertParagraph(line, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs right after the\n# \"GIS & Geospatial Analysis Consulting\" paragraph (under the Siege\n# Analytics / PARTNER heading), before the existing \"\u2022 Lead comprehensive\n# research initiatives...\" bullet.\n\n$d = $word.ActiveDocument\n\n$targetText = \"GIS & Geospatial Analysis Consulting\"\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($text -eq $targetText) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find paragraph: $targetText\"\n}\n\n$lines = @(\n    \"\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n    \"\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n    \"\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n)\n\n$insertAfterIndex = $targetIndex\nforeach ($line in $lines) {\n    $p = $d.Paragraphs($insertAfterIndex)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs($insertAfterIndex + 1)\n    $newPara.Range.Text = $line\n    $insertAfterIndex = $insertAfterIndex + 1\n}\n"}
